$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap B6/B7 values (role column)
$ws.Range("B6").Value = "tutor"
$ws.Range("B7").Value = "d"

# Update tutor email for row 6
$ws.Range("D6").Value = "tutor53@nkt.com"

# Update picture path for row 6
$ws.Range("K6").Value = "D:\test\picturesprofile1\bharat3.jpg"

# Update dates for row 6
$ws.Range("P6").Value = "07/12/2022"
$ws.Range("Q6").Value = "31/01/2023"

# Update sheet view (active selection / scroll position)
$excel.ActiveWindow.ScrollColumn = 4
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("L22").Select()

